$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 480.6
$ws.Range("I33").Value = 475.75
$ws.Range("K33").Value = 475.75
$ws.Range("M33").Value = -246.75
$ws.Range("H70").Value = 3879.3333
$ws.Range("I70").Value = 4091.4614
$ws.Range("J70").Value = 3628.6365
$ws.Range("K70").Value = 12274.3842
$ws.Range("L70").Value = 10885.9095
$ws.Range("M70").Value = -12004.3842
$ws.Range("N70").Value = -11425.9095
$ws.Range("H73").Value = 3879.3333
$ws.Range("I73").Value = 4091.4614
$ws.Range("J73").Value = 3628.6365
$ws.Range("K73").Value = 12274.3842
$ws.Range("L73").Value = 10885.9095
$ws.Range("M73").Value = -11338.3842
$ws.Range("N73").Value = -12757.9095
$ws.Range("H80").Value = 1688.2858
$ws.Range("I80").Value = 932.8333
$ws.Range("J80").Value = 2254.875
$ws.Range("K80").Value = 2798.4999
$ws.Range("L80").Value = 6764.625
$ws.Range("M80").Value = -1800.4999
$ws.Range("N80").Value = -8760.625
$ws.Range("H83").Value = 1688.2858
$ws.Range("I83").Value = 932.8333
$ws.Range("J83").Value = 2254.875
$ws.Range("K83").Value = 8395.4997
$ws.Range("L83").Value = 20293.875
$ws.Range("M83").Value = -3403.4997
$ws.Range("N83").Value = -30277.875
$ws.Range("H97").Value = 2202
$ws.Range("J97").Value = 2202
$ws.Range("L97").Value = 6606
$ws.Range("N97").Value = -7598
$ws.Range("H98").Value = 946
$ws.Range("I98").Value = 648.7646999999999
$ws.Range("J98").Value = 5999
$ws.Range("K98").Value = 648.7646999999999
$ws.Range("L98").Value = 5999
$ws.Range("M98").Value = 849.2353000000001
$ws.Range("N98").Value = -8995
$ws.Range("H100").Value = 3919
$ws.Range("I100").Value = 1796.2222
$ws.Range("K100").Value = 1796.2222
$ws.Range("M100").Value = -1255.2222
$ws.Range("H103").Value = 45455780
$ws.Range("I103").Value = 694
$ws.Range("J103").Value = 55556910
$ws.Range("K103").Value = 2082
$ws.Range("L103").Value = 166670730
$ws.Range("M103").Value = -1496
$ws.Range("N103").Value = -166671902
$ws.Range("H122").Value = 946
$ws.Range("I122").Value = 648.7646999999999
$ws.Range("J122").Value = 5999
$ws.Range("K122").Value = 1946.2941
$ws.Range("L122").Value = 17997
$ws.Range("M122").Value = 503.7059000000002
$ws.Range("N122").Value = -22897
$ws.Range("H134").Value = 113948.09
$ws.Range("J134").Value = 107393.1
$ws.Range("L134").Value = 107393.1
$ws.Range("N134").Value = -117533.1
$ws.Range("H138").Value = 2612.8333
$ws.Range("I138").Value = 1731.0625
$ws.Range("J138").Value = 3155.4614
$ws.Range("K138").Value = 5193.1875
$ws.Range("L138").Value = 9466.3842
$ws.Range("M138").Value = -53.1875
$ws.Range("N138").Value = -19746.3842

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1996.1143
$ws.Range("I2").Value = 1878.5
$ws.Range("K2").Value = 1878.5
$ws.Range("M2").Value = -1765.5
$ws.Range("H45").Value = 7147.76
$ws.Range("I45").Value = 12073.9
$ws.Range("J45").Value = 3863.6667
$ws.Range("K45").Value = 12073.9
$ws.Range("L45").Value = 3863.6667
$ws.Range("M45").Value = -11696.9
$ws.Range("N45").Value = -4617.6667
$ws.Range("H61").Value = 1823.375
$ws.Range("I61").Value = 1721.2354
$ws.Range("K61").Value = 1721.2354
$ws.Range("M61").Value = -1509.2354
$ws.Range("H63").Value = 2905
$ws.Range("I63").Value = 2359.2307
$ws.Range("K63").Value = 2359.2307
$ws.Range("M63").Value = -1673.2307
$ws.Range("H66").Value = 2905
$ws.Range("I66").Value = 2359.2307
$ws.Range("K66").Value = 11796.1535
$ws.Range("M66").Value = -8364.1535
$ws.Range("H74").Value = 2035.1111
$ws.Range("I74").Value = 1690.1666
$ws.Range("J74").Value = 2725
$ws.Range("K74").Value = 1690.1666
$ws.Range("L74").Value = 2725
$ws.Range("M74").Value = -816.1666
$ws.Range("N74").Value = -4473
$ws.Range("H77").Value = 2035.1111
$ws.Range("I77").Value = 1690.1666
$ws.Range("J77").Value = 2725
$ws.Range("K77").Value = 8450.833000000001
$ws.Range("L77").Value = 13625
$ws.Range("M77").Value = -4082.833000000001
$ws.Range("N77").Value = -22361
$ws.Range("H116").Value = 1996.1143
$ws.Range("I116").Value = 1878.5
$ws.Range("K116").Value = 1878.5
$ws.Range("M116").Value = 415.5
$ws.Range("H136").Value = 1823.375
$ws.Range("I136").Value = 1721.2354
$ws.Range("K136").Value = 5163.706200000001
$ws.Range("M136").Value = -2613.706200000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1996.1143
$ws.Range("I3").Value = 1878.5
$ws.Range("K3").Value = 1878.5
$ws.Range("M3").Value = -1764.5
$ws.Range("H20").Value = 5344.8057
$ws.Range("I20").Value = 7066.0835
$ws.Range("J20").Value = 1902.25
$ws.Range("K20").Value = 7066.0835
$ws.Range("L20").Value = 1902.25
$ws.Range("M20").Value = -6819.0835
$ws.Range("H107").Value = 114582.664
$ws.Range("I107").Value = 252962.5
$ws.Range("K107").Value = 252962.5
$ws.Range("M107").Value = -251042.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30820.656
$ws.Range("J31").Value = 3203.8
$ws.Range("L31").Value = 3203.8
$ws.Range("N31").Value = -3793.8
$ws.Range("H34").Value = 30820.656
$ws.Range("J34").Value = 3203.8
$ws.Range("L34").Value = 3203.8
$ws.Range("N34").Value = -3607.8
$ws.Range("H58").Value = 1482.625
$ws.Range("J58").Value = 1380.5
$ws.Range("L58").Value = 1380.5
$ws.Range("N58").Value = -1786.5
$ws.Range("H132").Value = 3602.926
$ws.Range("I132").Value = 3549.1924
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 10647.5772
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -8117.5772
$ws.Range("N132").Value = -20060
$ws.Range("H134").Value = 2323.375
$ws.Range("I134").Value = 1701.0731
$ws.Range("K134").Value = 5103.219300000001
$ws.Range("M134").Value = -2568.219300000001
$ws.Range("H136").Value = 1482.625
$ws.Range("J136").Value = 1380.5
$ws.Range("L136").Value = 4141.5
$ws.Range("N136").Value = -9241.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2233
$ws.Range("I3").Value = 2233
$ws.Range("K3").Value = 6699
$ws.Range("M3").Value = -6587
$ws.Range("H25").Value = 251499.75
$ws.Range("J25").Value = 2000
$ws.Range("L25").Value = 6000
$ws.Range("N25").Value = -6338
$ws.Range("H30").Value = 251499.75
$ws.Range("J30").Value = 2000
$ws.Range("L30").Value = 6000
$ws.Range("N30").Value = -6204
$ws.Range("H108").Value = 1405.4
$ws.Range("I108").Value = 1394.8889
$ws.Range("K108").Value = 4184.6667
$ws.Range("M108").Value = -1304.6667

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 8775
$ws.Range("I21").Value = 8633.333000000001
$ws.Range("K21").Value = 8633.333000000001
$ws.Range("M21").Value = -8460.333000000001
$ws.Range("H30").Value = 8775
$ws.Range("I30").Value = 8633.333000000001
$ws.Range("K30").Value = 8633.333000000001
$ws.Range("M30").Value = -8528.333000000001
$ws.Range("H122").Value = 3971.5
$ws.Range("I122").Value = 3962.1667
$ws.Range("J122").Value = 3999.5
$ws.Range("K122").Value = 11886.5001
$ws.Range("L122").Value = 11998.5
$ws.Range("M122").Value = -9436.500100000001
$ws.Range("N122").Value = -16898.5
$ws.Range("H126").Value = 2863.8333
$ws.Range("I126").Value = 3811
$ws.Range("J126").Value = 2390.25
$ws.Range("K126").Value = 11433
$ws.Range("L126").Value = 7170.75
$ws.Range("M126").Value = -8963
$ws.Range("N126").Value = -12110.75
$ws.Range("H132").Value = 3518.26
$ws.Range("I132").Value = 2643.0244
$ws.Range("J132").Value = 7505.4443
$ws.Range("K132").Value = 7929.073199999999
$ws.Range("L132").Value = 22516.3329
$ws.Range("M132").Value = -5399.073199999999
$ws.Range("N132").Value = -27576.3329
$ws.Range("H136").Value = 48435.6
$ws.Range("J136").Value = 48435.6
$ws.Range("L136").Value = 145306.8
$ws.Range("N136").Value = -150406.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14104.871
$ws.Range("I7").Value = 32025.1
$ws.Range("J7").Value = 5571.4287
$ws.Range("K7").Value = 32025.1
$ws.Range("L7").Value = 5571.4287
$ws.Range("M7").Value = -31913.1
$ws.Range("N7").Value = -5795.4287
$ws.Range("H122").Value = 147703.36
$ws.Range("I122").Value = 203851.6
$ws.Range("K122").Value = 611554.8
$ws.Range("M122").Value = -609104.8
$ws.Range("H126").Value = 14104.871
$ws.Range("I126").Value = 32025.1
$ws.Range("J126").Value = 5571.4287
$ws.Range("K126").Value = 96075.29999999999
$ws.Range("L126").Value = 16714.2861
$ws.Range("M126").Value = -93605.29999999999
$ws.Range("N126").Value = -21654.2861
$ws.Range("H132").Value = 3731.282
$ws.Range("I132").Value = 3101.9333
$ws.Range("J132").Value = 5829.1113
$ws.Range("K132").Value = 9305.7999
$ws.Range("L132").Value = 17487.3339
$ws.Range("M132").Value = -6775.7999
$ws.Range("N132").Value = -22547.3339

Write-Host "Applied all updates"